$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.373.43'
$ws.Range("E2").Value = '  +0.18%  '
$ws.Range("D3").Value = '1.879.89'
$ws.Range("E3").Value = '  +0.32%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.002'
$ws.Range("E4").Value = '  +0.20%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '0.7210'
$ws.Range("E5").Value = '  +1.50%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '242.69'
$ws.Range("E6").Value = '  +0.29%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '1.002'
$ws.Range("E7").Value = '  +0.20%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.08023'
$ws.Range("E8").Value = '  +2.66%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3145'
$ws.Range("E9").Value = '  +1.37%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '25.02'
$ws.Range("E10").Value = '  -0.49%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.08188'
$ws.Range("E11").Value = '  -2.70%  '
$ws.Range("D12").Value = '1.872.67'
$ws.Range("E12").Value = '  +0.48%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '94.80'
$ws.Range("E13").Value = '  +4.03%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.236'
$ws.Range("E14").Value = '  -0.09%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.7139'
$ws.Range("E15").Value = '  +0.26%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '6.416'
$ws.Range("E16").Value = '  +5.55%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.000008507'
$ws.Range("E17").Value = '  +2.02%  '
$ws.Range("D18").Value = '29.358.84'
$ws.Range("E18").Value = '  +0.13%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '244.52'
$ws.Range("E19").Value = '  +1.84%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.31'
$ws.Range("E20").Value = '  +0.62%  '
$ws.Range("E21").Value = '  +0.06%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '7.748'
$ws.Range("E22").Value = '  -0.07%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.003'
$ws.Range("E23").Value = '  +0.31%  '
$ws.Range("E24").Value = '  +0.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '9.047'
$ws.Range("E25").Value = '  +0.31%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '162.50'
$ws.Range("E26").Value = '  -0.03%  '
$ws.Range("E27").Value = '  +0.14%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '1.505'
$ws.Range("E28").Value = '  -0.04%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '4.410'
$ws.Range("E29").Value = '  +0.31%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '4.290'
$ws.Range("E30").Value = '  -1.02%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.226'
$ws.Range("E31").Value = '  -5.47%  '
$ws.Range("E32").Value = '  -0.02%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.942'
$ws.Range("E33").Value = '  -0.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.7681'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.179'
$ws.Range("E35").Value = '  +0.06%  '
$ws.Range("E36").Value = '  +0.73%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.01873'
$ws.Range("E37").Value = '  -0.14%  '
$ws.Range("D38").Value = '1.276.18'
$ws.Range("E38").Value = '  +4.34%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '2.753'
$ws.Range("E39").Value = '  +0.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '6.441'
$ws.Range("E40").Value = '  -0.72%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '113.20'
$ws.Range("E41").Value = '  +3.71%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.9109'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '74.40'
$ws.Range("E43").Value = '  +2.73%  '
$ws.Range("E44").Value = '  +7.01%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.002'
$ws.Range("E45").Value = '  +0.22%  '
$ws.Range("D46").Value = '2.020.67'
$ws.Range("E46").Value = '  -0.12%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.5230'
$ws.Range("E47").Value = '  +0.62%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.805'
$ws.Range("E48").Value = '  +0.43%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.506'
$ws.Range("E49").Value = '  +0.84%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.4347'
$ws.Range("E50").Value = '  +0.53%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '7.103'
$ws.Range("E51").Value = '  +0.35%  '
